$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1644736842105263
$ws.Range("C2").Value = 0.6085526315789473
$ws.Range("J2").Value = 0.006578947368421052
$ws.Range("P2").Value = 0.1217105263157895
$ws.Range("S2").Value = 0.09868421052631579
$ws.Range("C3").Value = 0.005319148936170213
$ws.Range("J3").Value = 0.03191489361702127
$ws.Range("P3").Value = 0.8191489361702128
$ws.Range("S3").Value = 0.1436170212765958
$ws.Range("J4").Value = 0.0392156862745098
$ws.Range("P4").Value = 0.7647058823529411
$ws.Range("S4").Value = 0.196078431372549
$ws.Range("B6").Value = 0.09056603773584905
$ws.Range("D6").Value = 0.01886792452830189
$ws.Range("F6").Value = 0.09811320754716982
$ws.Range("J6").Value = 0.2792452830188679
$ws.Range("O6").Value = 0.01132075471698113
$ws.Range("Q6").Value = 0.139622641509434
$ws.Range("R6").Value = 0.04150943396226415
$ws.Range("S6").Value = 0.3207547169811321
$ws.Range("B7").Value = 0.1294117647058824
$ws.Range("D7").Value = 0.02941176470588235
$ws.Range("F7").Value = 0.04117647058823529
$ws.Range("J7").Value = 0.1
$ws.Range("O7").Value = 0.01176470588235294
$ws.Range("Q7").Value = 0.2176470588235294
$ws.Range("R7").Value = 0.05882352941176471
$ws.Range("S7").Value = 0.4117647058823529
$ws.Range("B8").Value = 0.108256880733945
$ws.Range("D8").Value = 0.02201834862385321
$ws.Range("F8").Value = 0.06972477064220184
$ws.Range("J8").Value = 0.08256880733944955
$ws.Range("O8").Value = 0.02018348623853211
$ws.Range("Q8").Value = 0.1889908256880734
$ws.Range("R8").Value = 0.08073394495412844
$ws.Range("S8").Value = 0.4275229357798165
$ws.Range("B9").Value = 0.1032258064516129
$ws.Range("D9").Value = 0.02580645161290323
$ws.Range("F9").Value = 0.05161290322580645
$ws.Range("J9").Value = 0.09032258064516129
$ws.Range("O9").Value = 0.02580645161290323
$ws.Range("Q9").Value = 0.2258064516129032
$ws.Range("R9").Value = 0.09032258064516129
$ws.Range("S9").Value = 0.3870967741935484
$ws.Range("B10").Value = 0.1120917917034422
$ws.Range("D10").Value = 0.02383053839364519
$ws.Range("E10").Value = 0.00176522506619594
$ws.Range("F10").Value = 0.06619593998234775
$ws.Range("J10").Value = 0.08296557811120918
$ws.Range("O10").Value = 0.01588702559576346
$ws.Range("Q10").Value = 0.2391879964695499
$ws.Range("R10").Value = 0.0794351279788173
$ws.Range("S10").Value = 0.3786407766990291
$ws.Range("G11").Value = 0.1666666666666667
$ws.Range("J11").Value = 0.06944444444444445
$ws.Range("K11").Value = 0.2256944444444444
$ws.Range("L11").Value = 0.5277777777777778
$ws.Range("S11").Value = 0.01041666666666667
$ws.Range("G12").Value = 0.7077922077922078
$ws.Range("J12").Value = 0.2337662337662338
$ws.Range("L12").Value = 0.01948051948051948
$ws.Range("S12").Value = 0.03896103896103896
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.35
$ws.Range("S13").Value = 0.025
$ws.Range("F15").Value = 0.03791469194312796
$ws.Range("H15").Value = 0.1943127962085308
$ws.Range("I15").Value = 0.07109004739336493
$ws.Range("J15").Value = 0.3507109004739337
$ws.Range("K15").Value = 0.07582938388625593
$ws.Range("M15").Value = 0.02369668246445497
$ws.Range("O15").Value = 0.06635071090047394
$ws.Range("S15").Value = 0.1800947867298578
$ws.Range("F16").Value = 0.03111111111111111
$ws.Range("H16").Value = 0.2133333333333333
$ws.Range("I16").Value = 0.09777777777777778
$ws.Range("J16").Value = 0.3911111111111111
$ws.Range("K16").Value = 0.09333333333333334
$ws.Range("M16").Value = 0.01777777777777778
$ws.Range("O16").Value = 0.04888888888888889
$ws.Range("S16").Value = 0.1066666666666667
$ws.Range("F17").Value = 0.03909465020576132
$ws.Range("H17").Value = 0.2263374485596708
$ws.Range("I17").Value = 0.0720164609053498
$ws.Range("J17").Value = 0.4053497942386831
$ws.Range("K17").Value = 0.08847736625514403
$ws.Range("M17").Value = 0.00205761316872428
$ws.Range("O17").Value = 0.07613168724279835
$ws.Range("S17").Value = 0.09053497942386832
$ws.Range("F18").Value = 0.03550295857988166
$ws.Range("H18").Value = 0.2366863905325444
$ws.Range("I18").Value = 0.07100591715976332
$ws.Range("J18").Value = 0.3668639053254438
$ws.Range("K18").Value = 0.08284023668639054
$ws.Range("M18").Value = 0.01775147928994083
$ws.Range("O18").Value = 0.1005917159763314
$ws.Range("S18").Value = 0.08875739644970414
$ws.Range("F19").Value = 0.0347682119205298
$ws.Range("H19").Value = 0.2574503311258278
$ws.Range("I19").Value = 0.06043046357615894
$ws.Range("J19").Value = 0.3286423841059603
$ws.Range("K19").Value = 0.1043046357615894
$ws.Range("M19").Value = 0.02483443708609271
$ws.Range("O19").Value = 0.05877483443708609
$ws.Range("S19").Value = 0.130794701986755

Write-Output "Applied 105 cell updates"
